$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.620.37"
$ws.Cells.Item(2, 5).Value = "  +0.81%  "
$ws.Cells.Item(3, 4).Value = "3.382.92"
$ws.Cells.Item(3, 5).Value = "  -0.05%  "
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.16%  "
$ws.Cells.Item(5, 4).Value = "577.39"
$ws.Cells.Item(5, 5).Value = "  +0.86%  "
$ws.Cells.Item(6, 4).Value = "136.77"
$ws.Cells.Item(6, 5).Value = "  +0.37%  "
$ws.Cells.Item(7, 5).Value = "  -0.10%  "
$ws.Cells.Item(8, 4).Value = "3.382.12"
$ws.Cells.Item(8, 5).Value = "  -0.04%  "
$ws.Cells.Item(9, 5).Value = "  -0.68%  "
$ws.Cells.Item(10, 4).Value = "7.47"
$ws.Cells.Item(10, 5).Value = "  -1.47%  "
$ws.Cells.Item(11, 5).Value = "  +1.82%  "
$ws.Cells.Item(12, 5).Value = "  +0.61%  "
$ws.Cells.Item(13, 4).Value = "3.959.73"
$ws.Cells.Item(13, 5).Value = "  -0.26%  "
$ws.Cells.Item(14, 5).Value = "  +1.78%  "
$ws.Cells.Item(15, 5).Value = "  +1.14%  "
$ws.Cells.Item(16, 4).Value = "3.387.88"
$ws.Cells.Item(16, 5).Value = "  -0.23%  "
$ws.Cells.Item(17, 4).Value = "25.63"
$ws.Cells.Item(17, 5).Value = "  +1.84%  "
$ws.Cells.Item(18, 4).Value = "61.707.71"
$ws.Cells.Item(18, 5).Value = "  +0.52%  "
$ws.Cells.Item(19, 4).Value = "14.18"
$ws.Cells.Item(19, 5).Value = "  +1.07%  "
$ws.Cells.Item(20, 4).Value = "9.50"
$ws.Cells.Item(20, 5).Value = "  +0.90%  "
$ws.Cells.Item(21, 4).Value = "5.80"
$ws.Cells.Item(21, 5).Value = "  +0.12%  "
$ws.Cells.Item(22, 4).Value = "379.11"
$ws.Cells.Item(22, 5).Value = "  +1.26%  "
$ws.Cells.Item(23, 4).Value = "0.560"
$ws.Cells.Item(23, 5).Value = "  -1.40%  "
$ws.Cells.Item(24, 4).Value = "3.524.89"
$ws.Cells.Item(24, 5).Value = "  +0.01%  "
$ws.Cells.Item(25, 5).Value = "  -0.12%  "
$ws.Cells.Item(26, 4).Value = "71.22"
$ws.Cells.Item(26, 5).Value = "  +0.68%  "
$ws.Cells.Item(27, 4).Value = "0.0000125"
$ws.Cells.Item(27, 5).Value = "  +5.88%  "
$ws.Cells.Item(28, 5).Value = "  +4.20%  "
$ws.Cells.Item(29, 4).Value = "7.60"
$ws.Cells.Item(29, 5).Value = "  -1.65%  "
$ws.Cells.Item(30, 5).Value = "  -0.01%  "
$ws.Cells.Item(31, 4).Value = "8.18"
$ws.Cells.Item(31, 5).Value = "  +0.68%  "
$ws.Cells.Item(32, 5).Value = "  +2.33%  "
$ws.Cells.Item(33, 5).Value = "  +0.63%  "
$ws.Cells.Item(34, 5).Value = "  +0.04%  "
$ws.Cells.Item(35, 4).Value = "23.36"
$ws.Cells.Item(35, 5).Value = "  -0.32%  "
$ws.Cells.Item(36, 4).Value = "5.35"
$ws.Cells.Item(36, 5).Value = "  -4.03%  "
$ws.Cells.Item(37, 4).Value = "1.56"
$ws.Cells.Item(37, 5).Value = "  -0.86%  "
$ws.Cells.Item(38, 4).Value = "6.84"
$ws.Cells.Item(38, 5).Value = "  -1.36%  "
$ws.Cells.Item(39, 4).Value = "164.65"
$ws.Cells.Item(39, 5).Value = "  +0.97%  "
$ws.Cells.Item(40, 4).Value = "0.0785"
$ws.Cells.Item(40, 5).Value = "  -0.66%  "
$ws.Cells.Item(41, 2).Value = "Mantle"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(41, 4).Value = "0.781"
$ws.Cells.Item(41, 5).Value = "  +2.54%  "
$ws.Cells.Item(42, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(42, 4).Value = "0.999"
$ws.Cells.Item(42, 5).Value = "  -0.25%  "
$ws.Cells.Item(43, 5).Value = "  +1.97%  "
$ws.Cells.Item(44, 5).Value = "  +7.33%  "
$ws.Cells.Item(45, 2).Value = "Filecoin"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(45, 4).Value = "4.41"
$ws.Cells.Item(45, 5).Value = "  -0.25%  "
$ws.Cells.Item(46, 2).Value = "OKB"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(46, 4).Value = "41.27"
$ws.Cells.Item(46, 5).Value = "  -0.24%  "
$ws.Cells.Item(47, 2).Value = "EnergySwap"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(47, 4).Value = "24.71"
$ws.Cells.Item(47, 5).Value = "  +6.11%  "
$ws.Cells.Item(48, 4).Value = "6.86"
$ws.Cells.Item(48, 5).Value = "  -1.73%  "
$ws.Cells.Item(49, 4).Value = "22.81"
$ws.Cells.Item(49, 5).Value = "  -0.81%  "
$ws.Cells.Item(50, 4).Value = "2.335.28"
$ws.Cells.Item(50, 5).Value = "  +5.63%  "
$ws.Cells.Item(51, 5).Value = "  +1.35%  "
